$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20, pushing existing rows 20:91 down to 21:92
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new weekly price record
$ws.Cells.Item(20, 1).Value  = 7
$ws.Cells.Item(20, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value  = "Ñuble"
$ws.Cells.Item(20, 4).Value  = 45222
$ws.Cells.Item(20, 5).Value  = 16
$ws.Cells.Item(20, 6).Value  = 100112026
$ws.Cells.Item(20, 7).Value  = "Haba"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 60
$ws.Cells.Item(20, 11).Value = 15000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 15000
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(20, 16).Value = 600
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
